$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.2277013571343607
$ws.Range("D2").Value = 0.8219812918541203

# Row 3
$ws.Range("C3").Value = -0.3311860086160919
$ws.Range("D3").Value = 0.7436378936987875

# Row 4
$ws.Range("C4").Value = 2.464124884423317
$ws.Range("D4").Value = 0.02202237479343205

# Row 5
$ws.Range("C5").Value = 1.74357819331541
$ws.Range("D5").Value = 0.09519261137781454

# Row 6
$ws.Range("C6").Value = -0.4478740411138314
$ws.Range("D6").Value = 0.658623399355208

# Row 7
$ws.Range("C7").Value = 2.529241004611813
$ws.Range("D7").Value = 0.01910456335690136

# Row 8
$ws.Range("C8").Value = 1.419894817059028
$ws.Range("D8").Value = 0.1696545707763584

# Row 9
$ws.Range("C9").Value = 2.885916100916389
$ws.Range("D9").Value = 0.008578059008442107

# Row 10
$ws.Range("C10").Value = 3.029290782299383
$ws.Range("D10").Value = 0.006161268625352712
$ws.Range("G10").Value = "Sí"

# Row 11
$ws.Range("C11").Value = -1.206239800730076
$ws.Range("D11").Value = 0.2405372707273088
